# Updates the cryptos list with refreshed price/volume data (and a reorder
# of three rows) as scraped on Tue May 16 17:56:06 UTC 2023.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- Price column (D): many of these strings parse as plain numbers (e.g.
# "1.001"), which would make Excel silently store them as numeric cells and
# mint a new cell style. Force the Price column to the Text number format
# first so the values stay text exactly like the source data, then restore
# the original (default) cell style so no spurious style diff is introduced.
$dCells = @('D2', 'D3', 'D4', 'D5', 'D8', 'D9', 'D10', 'D12', 'D13', 'D14', 'D15', 'D16', 'D19', 'D21', 'D22', 'D24', 'D25', 'D26', 'D27', 'D28', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($cellRef in $dCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '27.080.55'
$ws.Range('D3').Value = '1.826.03'
$ws.Range('D4').Value = '1.001'
$ws.Range('D5').Value = '311.22'
$ws.Range('D8').Value = '0.3677'
$ws.Range('D9').Value = '0.07226'
$ws.Range('D10').Value = '0.8437'
$ws.Range('D12').Value = '1.819.55'
$ws.Range('D13').Value = '6.667'
$ws.Range('D14').Value = '0.07054'
$ws.Range('D15').Value = '5.292'
$ws.Range('D16').Value = '89.79'
$ws.Range('D19').Value = '1.000'
$ws.Range('D21').Value = '27.143.06'
$ws.Range('D22').Value = '5.145'
$ws.Range('D24').Value = '2.044.91'
$ws.Range('D25').Value = '1.986'
$ws.Range('D26').Value = '151.62'
$ws.Range('D27').Value = '2.245'
$ws.Range('D28').Value = '18.25'
$ws.Range('D29').Value = '5.271'
$ws.Range('D30').Value = '116.83'
$ws.Range('D31').Value = '0.08734'
$ws.Range('D32').Value = '1.179'
$ws.Range('D33').Value = '0.7376'
$ws.Range('D34').Value = '4.424'
$ws.Range('D35').Value = '2.882'
$ws.Range('D36').Value = '0.9999'
$ws.Range('D40').Value = '7.313'
$ws.Range('D41').Value = '2.877'
$ws.Range('D42').Value = '0.1686'
$ws.Range('D43').Value = '0.5069'
$ws.Range('D44').Value = '8.570'
$ws.Range('D45').Value = '10.59'
$ws.Range('D46').Value = '1.950'
$ws.Range('D47').Value = '106.06'
$ws.Range('D48').Value = '0.4728'
$ws.Range('D49').Value = '0.9995'
$ws.Range('D50').Value = '0.06327'
$ws.Range('D51').Value = '1.654'

# Restore default (unstyled) look - copy the style from a cell that was
# never touched and still carries the original default style.
foreach ($cellRef in $dCells) {
    $ws.Range($cellRef).Style = $ws.Range('D6').Style
}

# -- Coin / Link / Volume columns (B, C, E): plain text updates.
$ws.Range('E2').Value = '  -2.44%  '
$ws.Range('E3').Value = '  -1.36%  '
$ws.Range('E4').Value = '  -1.18%  '
$ws.Range('E5').Value = '  -2.37%  '
$ws.Range('E8').Value = '  -2.02%  '
$ws.Range('E9').Value = '  -1.61%  '
$ws.Range('E10').Value = '  -3.75%  '
$ws.Range('E11').Value = '  -3.53%  '
$ws.Range('E12').Value = '  -1.74%  '
$ws.Range('E13').Value = '  -1.28%  '
$ws.Range('E14').Value = '  -1.21%  '
$ws.Range('E15').Value = '  -2.72%  '
$ws.Range('E16').Value = '  +1.04%  '
$ws.Range('E17').Value = '  -1.20%  '
$ws.Range('E18').Value = '  -2.82%  '
$ws.Range('E19').Value = '  -1.04%  '
$ws.Range('E20').Value = '  -3.60%  '
$ws.Range('E21').Value = '  -2.24%  '
$ws.Range('E22').Value = '  -1.36%  '
$ws.Range('E24').Value = '  -1.68%  '
$ws.Range('E25').Value = '  +0.10%  '
$ws.Range('E26').Value = '  -2.40%  '
$ws.Range('E27').Value = '  +3.14%  '
$ws.Range('E28').Value = '  -1.83%  '
$ws.Range('E29').Value = '  -1.58%  '
$ws.Range('E30').Value = '  -1.81%  '
$ws.Range('E32').Value = '  -3.97%  '
$ws.Range('E33').Value = '  -5.12%  '
$ws.Range('E34').Value = '  -2.59%  '
$ws.Range('E35').Value = '  -1.47%  '
$ws.Range('E36').Value = '  -1.19%  '
$ws.Range('E37').Value = '  -3.46%  '
$ws.Range('E38').Value = '  -1.53%  '
$ws.Range('E39').Value = '  -1.99%  '
$ws.Range('E40').Value = '  +2.05%  '
$ws.Range('E41').Value = '  -0.86%  '
$ws.Range('E42').Value = '  -0.68%  '
$ws.Range('E43').Value = '  -1.30%  '
$ws.Range('E44').Value = '  -2.93%  '
$ws.Range('E45').Value = '  -1.45%  '
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('E46').Value = '  +5.73%  '
$ws.Range('B47').Value = 'Quant'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('E47').Value = '  -1.21%  '
$ws.Range('B48').Value = 'Decentraland'
$ws.Range('C48').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('E48').Value = '  -0.81%  '
$ws.Range('E49').Value = '  -1.28%  '
$ws.Range('E50').Value = '  -2.21%  '
$ws.Range('E51').Value = '  -2.13%  '
